$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 32, shifting existing rows 32-38 down to 33-39
$ws.Rows.Item(32).Insert()

# Fill in the new row 32 with the new weekly record
$ws.Cells.Item(32, 1).Value = 3
$ws.Cells.Item(32, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(32, 3).Value = "Coquimbo"
$ws.Cells.Item(32, 4).Value = 44524
$ws.Cells.Item(32, 5).Value = 5
$ws.Cells.Item(32, 6).Value = 100112022
$ws.Cells.Item(32, 7).Value = "Arveja Verde"
$ws.Cells.Item(32, 8).Value = "Perfection"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 65
$ws.Cells.Item(32, 11).Value = 16000
$ws.Cells.Item(32, 12).Value = 17000
$ws.Cells.Item(32, 13).Value = 16538
$ws.Cells.Item(32, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(32, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(32, 16).Value = 662
$ws.Cells.Item(32, 17).Value = 25
$ws.Cells.Item(32, 18).Value = "Hortaliza"

# Copy the date cell style (s="2") from the row below into the new row's D cell
$ws.Cells.Item(33, 4).Copy()
$ws.Cells.Item(32, 4).PasteSpecial(-4122)
